# Weekly re-aggregation of the "Achicoria" sheet: rows 2..31 keep their
# row-constant columns (A,B,C,E,F,G,H,I,N,Q,R) but the per-record columns
# (D Fecha, J Volumen, K Precio minimo, L Precio maximo, M Precio promedio
# ponderado, O Origen, P Precio $/Kg) get redistributed across rows
# according to the new weekly ordering. Concretely every new row's set of
# values is an exact copy of some other (old) row's values, i.e. the block
# is a row permutation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new_row -> old_row (source of the D/J/K/L/M/O/P values for that row)
$perm = @{
    2 = 8;  3 = 17; 4 = 5;  5 = 24; 6 = 30; 7 = 6;  8 = 22; 9 = 18; 10 = 16;
    11 = 23; 12 = 10; 13 = 25; 14 = 15; 15 = 31; 16 = 2; 17 = 7; 18 = 14;
    19 = 26; 20 = 21; 21 = 13; 22 = 20; 23 = 27; 24 = 29; 25 = 19; 26 = 4;
    27 = 28; 28 = 3; 29 = 9; 30 = 12; 31 = 11
}

$cols = @(4, 10, 11, 12, 13, 15, 16)   # D, J, K, L, M, O, P

# Snapshot every source value BEFORE any writes (it's a permutation, not a
# simple shift, so we must not read a cell after it has already been
# overwritten).
$snapshot = @{}
foreach ($r in 2..31) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the permuted values back.
foreach ($r in 2..31) {
    $src = $perm[$r]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c]
    }
}
